# Apply the changes described by the commit:
#   - rename "log_type" header to "reporting_type" (on the existing sheet)
#   - add a new worksheet "Sampling_Time_Unit_Digits" with the
#     "sampling delay" (Delta t) digit-distribution table

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Rename the "log_type" column header to "reporting_type" ---
$ws1.Range("B1").Value = "reporting_type"

# --- Add the new worksheet after the existing one ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "Sampling_Time_Unit_Digits"

# Row 1: header row ("reporting_type", 0, 1, ..., 9)
# Copy A1:K1 from sheet1 first so the new cells inherit the same
# bold/centered/bordered style already used for headers.
$ws1.Range("A1:K1").Copy($ws2.Range("A1"))
$ws2.Cells.Item(1, 1).Value = "reporting_type"
$headerVals = @(0, 1, 2, 3, 4, 5, 6, 7, 8, 9)
for ($col = 2; $col -le 11; $col++) {
    $ws2.Cells.Item(1, $col).Value = $headerVals[$col - 2]
}

# Row 2: Selfreport
$ws1.Range("A2").Copy($ws2.Range("A2"))
$ws2.Cells.Item(2, 1).Value = "Selfreport"
$row2vals = @(24, 9, 8, 7, 4, 21, 8, 9, 7, 4)
for ($col = 2; $col -le 11; $col++) {
    $ws2.Cells.Item(2, $col).Value = $row2vals[$col - 2]
}

# Row 3: App
$ws1.Range("A2").Copy($ws2.Range("A3"))
$ws2.Cells.Item(3, 1).Value = "App"
$row3vals = @(8, 13, 12, 10, 7, 9, 10, 12, 12, 8)
for ($col = 2; $col -le 11; $col++) {
    $ws2.Cells.Item(3, $col).Value = $row3vals[$col - 2]
}

# Keep the original sheet as the active/selected tab.
$ws1.Activate()
